# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 79
    5  = 2584
    7  = 131
    8  = 5
    9  = 1323
    11 = 55
    13 = 1172
    14 = 339
    16 = 31
    18 = 106
    19 = 69
    21 = 2383
    22 = 23
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
